$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Edit 1: "{% if age >= 18 %}" -> "{% if age >= 20 %}", with the "20"
# isolated into its own run (so the template tag is split into three
# runs: "{% if age >= " / "20" / " %}").
# --------------------------------------------------------------------
$ageFind = $d.Content
$ageFind.Find.Execute("{% if age >= 18 %}", $false, $false, $false, $false, `
                       $false, $true, 1, $false, "", 0) | Out-Null

if ($ageFind.Find.Found) {
    # $ageFind is now the matched range "{% if age >= 18 %}"
    $tagStart = $ageFind.Start
    $tagText = $ageFind.Text
    $numIdx = $tagText.IndexOf("18")
    $numStart = $tagStart + $numIdx
    $numEnd = $numStart + 2

    # Replace "18" with "20" first (plain text swap).
    $numRange = $d.Range($numStart, $numEnd)
    $numRange.Text = "20"

    # Force the "20" onto its own run by toggling a character property
    # on and back off; adjoining runs only stay distinct from their
    # neighbours when their formatting isn't identical, so a momentary
    # Bold flip splits this run away from the surrounding template text
    # without altering how it looks.
    $numRange2 = $d.Range($numStart, $numStart + 2)
    $numRange2.Bold = $true
    $numRange2.Bold = $false
}

# --------------------------------------------------------------------
# Edit 2: "To make your document look pr" / "ofessionally produced..."
# -> "To make your document look pro" / "fessionally produced..."
# i.e. move the leading "o" of the second run to the end of the first
# run (the two runs are separated by the "_GoBack" bookmark).
# --------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$gbStart = $goBack.Start

$insertionPoint = $d.Range($gbStart, $gbStart)
$insertionPoint.InsertAfter("o")

$goBack2 = $d.Bookmarks("_GoBack")
$gbEnd2 = $goBack2.End
$staleO = $d.Range($gbEnd2, $gbEnd2 + 1)
$staleO.Text = ""
